# Add crime-count data for 2023-10-05 to output/violent-crime-full-year.xlsx
# Only column values for 2023 (and the two affected 2022/2020 cells on the
# Citywide/By-Neighborhood summary sheets) change; every other cell, style and
# sheet is left untouched.
$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = 'Citywide Totals'; Cell = 'J2'; Value = 5846 },
    @{ Sheet = 'Citywide Totals'; Cell = 'J3'; Value = 6252 },
    @{ Sheet = 'Citywide Totals'; Cell = 'G4'; Value = 1473 },
    @{ Sheet = 'Citywide Totals'; Cell = 'I4'; Value = 1774 },
    @{ Sheet = 'Citywide Totals'; Cell = 'J4'; Value = 1352 },
    @{ Sheet = 'Citywide Totals'; Cell = 'J5'; Value = 476 },
    @{ Sheet = 'Citywide Totals'; Cell = 'J6'; Value = 7994 },
    @{ Sheet = 'Citywide Totals'; Cell = 'G7'; Value = 24698 },
    @{ Sheet = 'Citywide Totals'; Cell = 'I7'; Value = 26230 },
    @{ Sheet = 'Citywide Totals'; Cell = 'J7'; Value = 21920 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J4'; Value = 97 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J7'; Value = 649 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J8'; Value = 1374 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J11'; Value = 345 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J13'; Value = 27 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J15'; Value = 244 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J16'; Value = 85 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J18'; Value = 182 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J19'; Value = 648 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J20'; Value = 452 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J29'; Value = 1216 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J31'; Value = 198 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J33'; Value = 1017 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J36'; Value = 300 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J37'; Value = 674 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J41'; Value = 141 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J42'; Value = 923 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J43'; Value = 177 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J44'; Value = 166 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J45'; Value = 33 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J50'; Value = 129 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J52'; Value = 548 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J53'; Value = 304 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J54'; Value = 427 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J55'; Value = 301 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J59'; Value = 26 },
    @{ Sheet = 'By Neighborhood'; Cell = 'G63'; Value = 273 },
    @{ Sheet = 'By Neighborhood'; Cell = 'I63'; Value = 243 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J63'; Value = 77 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J65'; Value = 553 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J67'; Value = 828 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J69'; Value = 49 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J73'; Value = 211 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J75'; Value = 67 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J77'; Value = 163 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J79'; Value = 626 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J83'; Value = 446 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J84'; Value = 187 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J85'; Value = 906 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J86'; Value = 138 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J91'; Value = 250 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J92'; Value = 68 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J93'; Value = 97 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J97'; Value = 186 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J98'; Value = 159 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J99'; Value = 347 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J100'; Value = 41 },
    @{ Sheet = 'By Neighborhood'; Cell = 'G101'; Value = 24698 },
    @{ Sheet = 'By Neighborhood'; Cell = 'I101'; Value = 26230 },
    @{ Sheet = 'By Neighborhood'; Cell = 'J101'; Value = 21920 },
    @{ Sheet = 'Auburn Gresham'; Cell = 'J2'; Value = 199 },
    @{ Sheet = 'Auburn Gresham'; Cell = 'J6'; Value = 210 },
    @{ Sheet = 'Auburn Gresham'; Cell = 'J7'; Value = 649 },
    @{ Sheet = 'Belmont Cragin'; Cell = 'J2'; Value = 103 },
    @{ Sheet = 'Belmont Cragin'; Cell = 'J6'; Value = 144 },
    @{ Sheet = 'Belmont Cragin'; Cell = 'J7'; Value = 345 },
    @{ Sheet = 'South Shore'; Cell = 'J3'; Value = 326 },
    @{ Sheet = 'South Shore'; Cell = 'J6'; Value = 265 },
    @{ Sheet = 'South Shore'; Cell = 'J7'; Value = 906 },
    @{ Sheet = 'Little Village'; Cell = 'J3'; Value = 168 },
    @{ Sheet = 'Little Village'; Cell = 'J7'; Value = 548 },
    @{ Sheet = 'Norwood Park'; Cell = 'J3'; Value = 13 },
    @{ Sheet = 'Norwood Park'; Cell = 'J7'; Value = 49 },
    @{ Sheet = 'Logan Square'; Cell = 'J4'; Value = 10 },
    @{ Sheet = 'Logan Square'; Cell = 'J6'; Value = 199 },
    @{ Sheet = 'Logan Square'; Cell = 'J7'; Value = 304 },
    @{ Sheet = 'Austin'; Cell = 'J3'; Value = 421 },
    @{ Sheet = 'Austin'; Cell = 'J6'; Value = 465 },
    @{ Sheet = 'Austin'; Cell = 'J7'; Value = 1374 },
    @{ Sheet = 'South Chicago'; Cell = 'J2'; Value = 135 },
    @{ Sheet = 'South Chicago'; Cell = 'J3'; Value = 163 },
    @{ Sheet = 'South Chicago'; Cell = 'J6'; Value = 121 },
    @{ Sheet = 'South Chicago'; Cell = 'J7'; Value = 446 },
    @{ Sheet = 'Garfield Park'; Cell = 'J2'; Value = 244 },
    @{ Sheet = 'Garfield Park'; Cell = 'J3'; Value = 339 },
    @{ Sheet = 'Garfield Park'; Cell = 'J6'; Value = 350 },
    @{ Sheet = 'Garfield Park'; Cell = 'J7'; Value = 1017 },
    @{ Sheet = 'Grand Crossing'; Cell = 'J2'; Value = 199 },
    @{ Sheet = 'Grand Crossing'; Cell = 'J3'; Value = 229 },
    @{ Sheet = 'Grand Crossing'; Cell = 'J6'; Value = 196 },
    @{ Sheet = 'Grand Crossing'; Cell = 'J7'; Value = 674 },
    @{ Sheet = 'New City'; Cell = 'J2'; Value = 162 },
    @{ Sheet = 'New City'; Cell = 'J6'; Value = 196 },
    @{ Sheet = 'New City'; Cell = 'J7'; Value = 553 },
    @{ Sheet = 'Woodlawn'; Cell = 'J2'; Value = 98 },
    @{ Sheet = 'Woodlawn'; Cell = 'J3'; Value = 136 },
    @{ Sheet = 'Woodlawn'; Cell = 'J4'; Value = 15 },
    @{ Sheet = 'Woodlawn'; Cell = 'J7'; Value = 347 },
    @{ Sheet = 'Gage Park'; Cell = 'J6'; Value = 54 },
    @{ Sheet = 'Gage Park'; Cell = 'J7'; Value = 198 },
    @{ Sheet = 'North Lawndale'; Cell = 'J2'; Value = 206 },
    @{ Sheet = 'North Lawndale'; Cell = 'J3'; Value = 313 },
    @{ Sheet = 'North Lawndale'; Cell = 'J6'; Value = 222 },
    @{ Sheet = 'North Lawndale'; Cell = 'J7'; Value = 828 },
    @{ Sheet = 'South Deering'; Cell = 'J3'; Value = 61 },
    @{ Sheet = 'South Deering'; Cell = 'J6'; Value = 56 },
    @{ Sheet = 'South Deering'; Cell = 'J7'; Value = 187 },
    @{ Sheet = 'Loop'; Cell = 'J5'; Value = 4 },
    @{ Sheet = 'Loop'; Cell = 'J6'; Value = 206 },
    @{ Sheet = 'Loop'; Cell = 'J7'; Value = 427 },
    @{ Sheet = 'Englewood'; Cell = 'J2'; Value = 369 },
    @{ Sheet = 'Englewood'; Cell = 'J4'; Value = 67 },
    @{ Sheet = 'Englewood'; Cell = 'J7'; Value = 1216 },
    @{ Sheet = 'Chatham'; Cell = 'J3'; Value = 190 },
    @{ Sheet = 'Chatham'; Cell = 'J6'; Value = 248 },
    @{ Sheet = 'Chatham'; Cell = 'J7'; Value = 648 },
    @{ Sheet = 'Irving Park'; Cell = 'J2'; Value = 51 },
    @{ Sheet = 'Irving Park'; Cell = 'J7'; Value = 166 },
    @{ Sheet = 'Hermosa'; Cell = 'J6'; Value = 78 },
    @{ Sheet = 'Hermosa'; Cell = 'J7'; Value = 141 },
    @{ Sheet = 'Humboldt Park'; Cell = 'J2'; Value = 199 },
    @{ Sheet = 'Humboldt Park'; Cell = 'J3'; Value = 189 },
    @{ Sheet = 'Humboldt Park'; Cell = 'J6'; Value = 477 },
    @{ Sheet = 'Humboldt Park'; Cell = 'J7'; Value = 923 },
    @{ Sheet = 'Boystown'; Cell = 'J5'; Value = 13 },
    @{ Sheet = 'Boystown'; Cell = 'J6'; Value = 27 },
    @{ Sheet = 'Lower West Side'; Cell = 'J3'; Value = 66 },
    @{ Sheet = 'Lower West Side'; Cell = 'J6'; Value = 154 },
    @{ Sheet = 'Lower West Side'; Cell = 'J7'; Value = 301 },
    @{ Sheet = 'Washington Park'; Cell = 'J3'; Value = 104 },
    @{ Sheet = 'Washington Park'; Cell = 'J6'; Value = 59 },
    @{ Sheet = 'Washington Park'; Cell = 'J7'; Value = 250 },
    @{ Sheet = 'Roseland'; Cell = 'J3'; Value = 214 },
    @{ Sheet = 'Roseland'; Cell = 'J4'; Value = 37 },
    @{ Sheet = 'Roseland'; Cell = 'J6'; Value = 182 },
    @{ Sheet = 'Roseland'; Cell = 'J7'; Value = 626 },
    @{ Sheet = 'Chicago Lawn'; Cell = 'J3'; Value = 156 },
    @{ Sheet = 'Chicago Lawn'; Cell = 'J6'; Value = 122 },
    @{ Sheet = 'Chicago Lawn'; Cell = 'J7'; Value = 452 },
    @{ Sheet = 'Calumet Heights'; Cell = 'J3'; Value = 39 },
    @{ Sheet = 'Calumet Heights'; Cell = 'J7'; Value = 182 },
    @{ Sheet = 'Grand Boulevard'; Cell = 'J2'; Value = 97 },
    @{ Sheet = 'Grand Boulevard'; Cell = 'J3'; Value = 97 },
    @{ Sheet = 'Grand Boulevard'; Cell = 'J6'; Value = 91 },
    @{ Sheet = 'Grand Boulevard'; Cell = 'J7'; Value = 300 },
    @{ Sheet = 'West Lawn'; Cell = 'J4'; Value = 9 },
    @{ Sheet = 'West Lawn'; Cell = 'J7'; Value = 97 },
    @{ Sheet = 'Wrigleyville'; Cell = 'J6'; Value = 21 },
    @{ Sheet = 'Wrigleyville'; Cell = 'J7'; Value = 41 },
    @{ Sheet = 'Brighton Park'; Cell = 'J2'; Value = 71 },
    @{ Sheet = 'Brighton Park'; Cell = 'J6'; Value = 102 },
    @{ Sheet = 'Brighton Park'; Cell = 'J7'; Value = 244 },
    @{ Sheet = 'Wicker Park'; Cell = 'J4'; Value = 9 },
    @{ Sheet = 'Wicker Park'; Cell = 'J6'; Value = 98 },
    @{ Sheet = 'Wicker Park'; Cell = 'J7'; Value = 159 },
    @{ Sheet = 'Lincoln Square'; Cell = 'J3'; Value = 34 },
    @{ Sheet = 'Lincoln Square'; Cell = 'J7'; Value = 129 },
    @{ Sheet = 'Portage Park'; Cell = 'J3'; Value = 53 },
    @{ Sheet = 'Portage Park'; Cell = 'J4'; Value = 14 },
    @{ Sheet = 'Portage Park'; Cell = 'J7'; Value = 211 },
    @{ Sheet = 'Montclare'; Cell = 'J3'; Value = 4 },
    @{ Sheet = 'Montclare'; Cell = 'J7'; Value = 26 },
    @{ Sheet = 'West Town'; Cell = 'J6'; Value = 129 },
    @{ Sheet = 'West Town'; Cell = 'J7'; Value = 186 },
    @{ Sheet = 'West Elsdon'; Cell = 'J2'; Value = 18 },
    @{ Sheet = 'West Elsdon'; Cell = 'J7'; Value = 68 },
    @{ Sheet = 'Streeterville'; Cell = 'J3'; Value = 21 },
    @{ Sheet = 'Streeterville'; Cell = 'J7'; Value = 138 },
    @{ Sheet = 'Pullman'; Cell = 'J2'; Value = 29 },
    @{ Sheet = 'Pullman'; Cell = 'J3'; Value = 20 },
    @{ Sheet = 'Pullman'; Cell = 'J7'; Value = 67 },
    @{ Sheet = 'Hyde Park'; Cell = 'J4'; Value = 19 },
    @{ Sheet = 'Hyde Park'; Cell = 'J7'; Value = 177 },
    @{ Sheet = 'Riverdale'; Cell = 'J4'; Value = 15 },
    @{ Sheet = 'Riverdale'; Cell = 'J7'; Value = 163 },
    @{ Sheet = 'Jackson Park'; Cell = 'J2'; Value = 12 },
    @{ Sheet = 'Jackson Park'; Cell = 'J7'; Value = 33 },
    @{ Sheet = 'Archer Heights'; Cell = 'J3'; Value = 23 },
    @{ Sheet = 'Archer Heights'; Cell = 'J7'; Value = 97 },
    @{ Sheet = 'Bucktown'; Cell = 'J6'; Value = 67 },
    @{ Sheet = 'Bucktown'; Cell = 'J7'; Value = 85 }
)

foreach ($chg in $changes) {
    $ws = $wb.Worksheets.Item($chg.Sheet)
    $ws.Range($chg.Cell).Value = $chg.Value
}
